# Generate Report for Handoff
# - Update status text "In Translation" -> "Ready for handoff"
# - Update the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# - Widen the datetime columns slightly (E/F on Overview, C on zh-cn/de-de)

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-06 03:06:13"

# Widen columns E and F (status/date columns) to match the new, wider content.
# The ColumnWidth property is quantized to whole pixels by this runtime, so we
# pick the closest representable value to the intended width.
$wsOverview.Columns.Item(5).ColumnWidth = 98 / 6
$wsOverview.Columns.Item(6).ColumnWidth = 98 / 6

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-06 03:06:06"
$wsZhCn.Columns.Item(3).ColumnWidth = 98 / 6

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-06 03:06:13"
$wsDeDe.Columns.Item(3).ColumnWidth = 98 / 6
